# #332 Fix cross references
# Update General sheet: rename product references from 'HtmlSC' to 'HSC'
$wb = $excel.ActiveWorkbook

$general = $wb.Worksheets.Item("General")
$general.Range("C2").Value = 'HSC shall read a single (configurable) HTML file'
$general.Range("C3").Value = 'HSC can be run as {gradle-url}-plugin.'
$general.Range("C4").Value = 'HSC can be called from the command line with arguments and options'

# Update Required Checks sheet: point cross references at the new
# {xrefAlgorithm...Checker} attribute syntax instead of <<...Checker>> macros
$required = $wb.Worksheets.Item("Required Checks")
$required.Range("C2").Value = 'Check all image tags if the referenced image files exist. See {xrefAlgorithmMissingImageFilesChecker}'
$required.Range("C3").Value = 'Check all internal links from anchor-tags (href="#XYZ") if the link targets "XYZ" are defined. See {xrefAlgorithmBrokenCrossReferencesChecker}'
$required.Range("C4").Value = 'either other html-files, pdf’s or similar. See {xrefAlgorithmMissingLocalResourcesChecker}'
$required.Range("C5").Value = 'Check all bookmark definitions (…​ id="XYZ") whether the id’s ("XYZ") are unique. See {xrefAlgorithmDuplicateIdChecker}'
$required.Range("C7").Value = 'in image-tags. See {xrefAlgorithmMissingImgAltAttributeChecker}'

# Restore view state: General sheet selection moves to column C,
# Quality Goals scrolls to C34, and Required Checks becomes the active tab
$general.Range("C1:C1048576").Select()

$qualityGoals = $wb.Worksheets.Item("Quality Goals")
$qualityGoals.Activate()
$qualityGoals.Range("C34").Select()

$required.Activate()
$required.Range("C9").Select()

Write-Output "done"
